$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B40").Value = "xxxxx"
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = "magnet for lid sensor"

$ws.Range("F41").Select()
